$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 339, pushing the existing rows
# (old 339..351) down to (341..353).
$ws.Rows.Item(339).Insert()
$ws.Rows.Item(339).Insert()

# Populate the first new row (339)
$ws.Cells.Item(339, 1).Value = 7
$ws.Cells.Item(339, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(339, 3).Value = "Ñuble"
$ws.Cells.Item(339, 4).Value = 44939
$ws.Cells.Item(339, 5).Value = 16
$ws.Cells.Item(339, 6).Value = 100112009
$ws.Cells.Item(339, 7).Value = "Acelga"
$ws.Cells.Item(339, 8).Value = "Sin especificar"
$ws.Cells.Item(339, 9).Value = "Primera"
$ws.Cells.Item(339, 10).Value = 300
$ws.Cells.Item(339, 11).Value = 600
$ws.Cells.Item(339, 12).Value = 650
$ws.Cells.Item(339, 13).Value = 625
$ws.Cells.Item(339, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(339, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(339, 16).Value = 625
$ws.Cells.Item(339, 17).Value = 1
$ws.Cells.Item(339, 18).Value = "Hortaliza"

# Populate the second new row (340)
$ws.Cells.Item(340, 1).Value = 7
$ws.Cells.Item(340, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(340, 3).Value = "Ñuble"
$ws.Cells.Item(340, 4).Value = 44939
$ws.Cells.Item(340, 5).Value = 16
$ws.Cells.Item(340, 6).Value = 100112009
$ws.Cells.Item(340, 7).Value = "Acelga"
$ws.Cells.Item(340, 8).Value = "Sin especificar"
$ws.Cells.Item(340, 9).Value = "Segunda"
$ws.Cells.Item(340, 10).Value = 200
$ws.Cells.Item(340, 11).Value = 500
$ws.Cells.Item(340, 12).Value = 500
$ws.Cells.Item(340, 13).Value = 500
$ws.Cells.Item(340, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(340, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(340, 16).Value = 500
$ws.Cells.Item(340, 17).Value = 1
$ws.Cells.Item(340, 18).Value = "Hortaliza"
